$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new work-log entries (rows 63 and 64) ---
# Copy formatting from the two rows immediately above so the new rows
# pick up the same date / number / wrapped-text styles already used
# throughout the table (style ids 3, 4, 9).
$ws.Range("E61:G61").Copy() | Out-Null
$ws.Range("E63").PasteSpecial(-4122) | Out-Null
$ws.Range("E62:G62").Copy() | Out-Null
$ws.Range("E64").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row 63: 05.01.2020 - 6h - resource bundles / preferences work
$ws.Cells.Item(63, 5).Value = 43835
$ws.Cells.Item(63, 6).Value = 6
$ws.Cells.Item(63, 7).Value = "Implementierung von Resource Bundles in die Applikation abgeschlossen`n.properties files mit weiteren Einträgen ergänzt`nBug behoben, bei dem sich die Position der Menüeinträge beim Ändern der Sprache geändert hat`nSpeicherung der Daten in java preferences umgesetzt`nAnpassungen an UI vorgenommen"

# Row 64: 06.01.2020 - 2h - radio buttons instead of checkboxes
$ws.Cells.Item(64, 5).Value = 43836
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(64, 7).Value = "Ersetzen von 2 checkboxes durch radio buttons`nAnpassungen an UI vorgenommen"

# Match the row heights Excel would have auto-fit for the wrapped text
# (5 lines / 2 lines at the default 14.4pt row height).
$ws.Rows.Item(63).RowHeight = 72
$ws.Rows.Item(64).RowHeight = 28.8

# --- Update the view state to match where the user ended up scrolled/selected ---
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("G64:G65").Select() | Out-Null
$ws.Range("G65").Activate() | Out-Null
